$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = 'در یک مطالعه اپیدمیولوژیک، اولین موردی که توجه محقق را جلب می‌کند چیست؟'
$ws.Rows.Item(2).AutoFit()
$ws.Range("A3").Value = 'کدام یک از ویژگی‌های زیر پتیدین را از مورفین متمایز می‌کند؟'
$ws.Rows.Item(3).AutoFit()
$ws.Range("A4").Value = 'کدام گزینه در مورد سیستادنوما نادرست است؟'
$ws.Rows.Item(4).AutoFit()
$ws.Range("A5").Value = 'یک التهاب ملتحمه دوطرفه عودکننده که با شروع هوای گرم در پسران جوان همراه با علائم سوزش، خارش و اشکریزش همراه با نواحی چندضلعی برجسته در ملتحمه پلکی رخ میدهد، کدام است؟'
$ws.Rows.Item(5).AutoFit()
$ws.Range("A6").Value = 'یک بیمار 50 ساله با سابقه طولانی مدت بیماری مزمن انسدادی ریه، به تدریج دچار درد در اندام‌های انتهایی، به‌ویژه در مچ‌های هر دو دست می‌شود. کاهش وزن 9 کیلوگرمی مشاهده می‌شود. پوست روی مچ‌ها گرم و اریتماتو است. کلابینگ دوطرفه وجود دارد. تصویربرداری ساده به‌عنوان ضخیم‌شدن پریوست و احتمال استئومیلیت گزارش شده است. شما باید:'
$ws.Rows.Item(6).AutoFit()
$ws.Range("A7").Value = 'هنگام انتخاب رنگ سیمان برای سیمان کردن یک روکش چینی، پودری که آزمایش می‌شود بهتر است با چه چیزی مخلوط شود؟'
$ws.Rows.Item(7).AutoFit()
$ws.Range("A8").Value = 'بی‌حسی گونه پس از شکستگی کمپلکس زیگوماتیک به چه علت است؟'
$ws.Rows.Item(8).AutoFit()
$ws.Range("A9").Value = 'بازقطبی شدن پتانسیل عمل عصبی به دلیل کدام مورد است؟'
$ws.Rows.Item(9).AutoFit()
$ws.Range("A10").Value = 'مردی با ضایعات تاولی روی صورت، زیر بغل و قفسه سینه مراجعه می‌کند. ابتدا این ضایعات فقط در دهان او وجود داشتند. این ضایعات گرد و بیضی شکل با مایع سروزی هستند و برخی از آنها "شل" هستند. هنگامی که فشار به ضایعه وارد می‌شود، مایع به صورت جانبی پخش می‌شود. تشخیص بالینی پمفیگوس ولگاریس (pv) داده شده است. کدام یک از گزینه‌های زیر سن معمول شروع این بیماری است؟'
$ws.Rows.Item(10).AutoFit()
$ws.Range("A11").Value = 'حداقل دوز اتینیل استرادیول در قرص های ضد بارداری خوراکی برای اثر ضد بارداری چقدر است؟'
$ws.Rows.Item(11).AutoFit()
$ws.Range("A12").Value = 'عملکرد فضا نگهدارها چیست؟'
$ws.Rows.Item(12).AutoFit()
$ws.Range("A13").Value = 'عصب تأمین کننده حسی ویژه ناحیه مشخص شده کدام است؟'
$ws.Rows.Item(13).AutoFit()
$ws.Range("A14").Value = 'جذب گلوکز توسط کبد پس از مصرف یک وعده غذایی حاوی کربوهیدرات افزایش می یابد زیرا'
$ws.Rows.Item(14).AutoFit()
$ws.Range("A15").Value = 'کدام گزینه در مورد پروپوفول نادرست است؟'
$ws.Rows.Item(15).AutoFit()
$ws.Range("A16").Value = 'در فیبروز کیستیک، سودوموناس آئروژینوزا شایع‌ترین ارگانیسم است. بعد از آن شایع‌ترین ارگانیسم کدام است؟'
$ws.Rows.Item(16).AutoFit()
$ws.Range("A17").Value = 'در طی کالبدشکافی، بهترین محل برای اندازه‌گیری دمای مرکزی بدن در مورد لواط کجاست؟'
$ws.Rows.Item(17).AutoFit()
$ws.Range("A18").Value = 'سلول‌های اپی‌تلیوئید در تمام موارد زیر دیده می‌شوند به جز'
$ws.Rows.Item(18).AutoFit()
$ws.Range("A19").Value = 'تنگی قفسه سینه در روزهای دوشنبه مشخصه کدام بیماری است؟'
$ws.Rows.Item(19).AutoFit()
$ws.Range("A20").Value = 'فیناستراید یک -'
$ws.Rows.Item(20).AutoFit()
$ws.Range("A21").Value = 'در مورد سرطان مری کدام گزینه صحیح است؟'
$ws.Rows.Item(21).AutoFit()
$ws.Range("A22").Value = 'آنتی بادی ضد هسته ای (ana) در تمام موارد زیر دیده می شود به جز-'
$ws.Rows.Item(22).AutoFit()
$ws.Range("A23").Value = 'اصول درمان عفونت چرکی سیستم اسکلتی عبارتند از'
$ws.Rows.Item(23).AutoFit()
$ws.Range("A24").Value = 'سلول اصلی مسئول انقباض اسکار کدام است؟'
$ws.Rows.Item(24).AutoFit()
$ws.Range("A25").Value = 'یک پسر ۱۰ ساله با انسداد بینی و خونریزی شدید متناوب از بینی مراجعه می‌کند. او یک توده سفت و صورتی رنگ در نازوفارنکس دارد. تمام موارد زیر به عنوان بررسی‌های تشخیصی در این مورد انجام می‌شوند به جز-'
$ws.Rows.Item(25).AutoFit()
$ws.Range("A26").Value = 'سوء مصرف کوکائین بیشتر شبیه به سوء مصرف کدام ماده است؟'
$ws.Rows.Item(26).AutoFit()
$ws.Range("A27").Value = 'در مورد مایکوپلاسما همه موارد زیر صحیح است، به جز:'
$ws.Rows.Item(27).AutoFit()
$ws.Range("A28").Value = 'کدام یک از علائم بلاک گانگلیون ستاره‌ای نیست؟'
$ws.Rows.Item(28).AutoFit()
$ws.Range("A29").Value = 'عبارت نادرست در مورد کولیت اولسراتیو کدام است؟'
$ws.Rows.Item(29).AutoFit()
$ws.Range("A30").Value = 'کدام مورد در مورد آب مروارید سرخجه صحیح نیست؟'
$ws.Rows.Item(30).AutoFit()
$ws.Range("A31").Value = 'شاخص توده بدنی همچنین به عنوان چه چیزی شناخته می شود؟'
$ws.Rows.Item(31).AutoFit()
$ws.Range("A32").Value = 'یک زن ۴۲ ساله در طول ۸ ماه گذشته احساس ناراحتی در ناحیه شکم داشته است. سی‌تی‌اسکن شکم یک توده کیستی ۶ سانتی‌متری در تخمدان راست با مناطق کوچک کلسیفیکاسیون نشان می‌دهد. پس از برداشتن تخمدان، توده در برش حاوی مو و سبوم فراوان است. در بررسی میکروسکوپی، توده دارای فضاهای غده‌ای پوشیده از اپیتلیوم ستونی، اپیتلیوم سنگفرشی با فولیکول‌های مو، غضروف و بافت همبند متراکم است. این توده از چه نوعی است؟'
$ws.Rows.Item(32).AutoFit()
$ws.Range("A33").Value = 'یک زن ۲۰ ساله با کم خونی، زردی خفیف به مدت ۲ سال، اسمیر محیطی نشان دهنده اسفروسیت ها، بهترین آزمایش برای انجام کدام است؟'
$ws.Rows.Item(33).AutoFit()
$ws.Range("A34").Value = 'کدام یک از موارد زیر نقص اولیه در سندرم پیر رابین است؟'
$ws.Rows.Item(34).AutoFit()
$ws.Range("A35").Value = 'کدام یک از ساختارهای جنینی زیر منشأ غده فوق کلیوی است؟'
$ws.Rows.Item(35).AutoFit()
$ws.Range("A36").Value = 'یک دختر ۱۵ ساله در طول ۱۰ سال گذشته چندین ندول روی پوست خود ایجاد کرده است. در معاینه فیزیکی، ۲۰ ندول سفت به اندازه ۰.۳ تا ۱ سانتیمتر به صورت پراکنده روی تنه و اندام‌های بیمار مشاهده می‌شود. همچنین ۱۲ ماکول قهوه‌ای روشن با میانگین قطر ۲ تا ۵ سانتیمتر روی پوست تنه وجود دارد. معاینه با slit-lamp ندول‌های رنگدانه‌دار در عنبیه را نشان می‌دهد. یک خواهر و برادر و یکی از والدین نیز به همین صورت درگیر هستند. تحلیل ژنتیکی یک جهش loss-of-function را نشان می‌دهد. کدام یک از الگوهای توارث زیر به احتمال زیاد در این خانواده وجود دارد؟'
$ws.Rows.Item(36).AutoFit()
$ws.Range("A37").Value = 'چاقی در تمام موارد زیر دیده می‌شود به جز؟'
$ws.Rows.Item(37).AutoFit()
$ws.Range("A38").Value = 'تمامی عبارات زیر در مورد گردش خون برونشی صحیح هستند، به جز:'
$ws.Rows.Item(38).AutoFit()
$ws.Range("A39").Value = 'شایع‌ترین محل آنوریسم مغزی کجاست؟'
$ws.Rows.Item(39).AutoFit()
$ws.Range("A40").Value = 'شایع‌ترین علت کم‌کاری پاراتیروئید چیست؟'
$ws.Rows.Item(40).AutoFit()
$ws.Range("A41").Value = 'چه زمانی در احیای نوزاد، ماساژ قفسه سینه را آغاز می‌کنید؟'
$ws.Rows.Item(41).AutoFit()
$ws.Range("A42").Value = 'ماسراسیون در مرگ کدام یک دیده می‌شود؟'
$ws.Rows.Item(42).AutoFit()
$ws.Range("A43").Value = 'داروی انتخابی برای آسم شدید حاد کدام است؟'
$ws.Rows.Item(43).AutoFit()
$ws.Range("A44").Value = 'یک فرد معتاد به مواد مخدر دچار شبه‌آنوریسم شده است. کدام یک از موارد زیر باید در روش‌های درمانی گنجانده شود؟'
$ws.Rows.Item(44).AutoFit()
$ws.Range("A45").Value = 'ضایعه پوست پیازی در عروق در کدام مورد دیده می‌شود؟'
$ws.Rows.Item(45).AutoFit()
$ws.Range("A46").Value = 'غلظت کانال سدیم در کدام ناحیه بیشترین مقدار را دارد؟'
$ws.Rows.Item(46).AutoFit()
$ws.Range("A47").Value = 'در طبقه‌بندی مور از اومفالوسل (اکزامفالوس)، نقص نوع i ناف کمتر از ________ سانتی‌متر است.'
$ws.Rows.Item(47).AutoFit()
$ws.Range("A48").Value = 'کدام یک از موارد زیر جزء دیافراگم شکاف صافی در گلومرول نیست؟'
$ws.Rows.Item(48).AutoFit()
$ws.Range("A49").Value = 'در کدام یک از شرایط زیر داکتیلیت دیده نمی‌شود؟'
$ws.Rows.Item(49).AutoFit()
$ws.Range("A50").Value = 'علت ترمیم دندان با کانتور بیش از حد چیست؟'
$ws.Rows.Item(50).AutoFit()
$ws.Range("A51").Value = 'نادرست در مورد دستورالعمل‌های چاه بهداشتی'
$ws.Rows.Item(51).AutoFit()
$ws.Range("A52").Value = 'نورپلنت حاوی چند کپسول لوونورژسترل است:'
$ws.Rows.Item(52).AutoFit()
$ws.Range("A53").Value = 'در مورد میازیس بینی کدام گزینه صحیح است؟'
$ws.Rows.Item(53).AutoFit()
$ws.Range("A54").Value = 'هسته رافه در ساقه مغز مسئول ترشح کدام یک از مواد شیمیایی زیر است؟'
$ws.Rows.Item(54).AutoFit()
$ws.Range("A55").Value = 'همه موارد زیر از ویژگی‌های بیماری اوریون هستند به جز _______'
$ws.Rows.Item(55).AutoFit()
$ws.Range("A56").Value = 'کدام اختلال با تشنج، آب مروارید و عقب ماندگی ذهنی همراه است؟'
$ws.Rows.Item(56).AutoFit()
$ws.Range("A57").Value = 'کدام یک از موارد زیر در لوکوپوئزیس مشاهده نمی‌شود؟'
$ws.Rows.Item(57).AutoFit()
$ws.Range("A58").Value = 'مرواریدهای عنبیه در یوئیت ناشی از کدام یک دیده می‌شوند؟'
$ws.Rows.Item(58).AutoFit()
$ws.Range("A59").Value = 'در یک بیمار اخیراً تشخیص داده شده با دیابت نوع 1، معاینه ته چشم در چه زمانی انجام می‌شود؟'
$ws.Rows.Item(59).AutoFit()
$ws.Range("A60").Value = 'در عمل فاترگیل موارد زیر انجام می‌شوند به جز:'
$ws.Rows.Item(60).AutoFit()
$ws.Range("A61").Value = 'همه موارد زیر از ویژگی‌های کیست کاذب لوزالمعده هستند به جز:'
$ws.Rows.Item(61).AutoFit()
$ws.Range("A62").Value = 'کدام یک از عبارات زیر در مورد آزمایش تشخیصی oncotype نادرست نیست؟'
$ws.Rows.Item(62).AutoFit()
$ws.Range("A63").Value = 'آرایش حروف چینی در کدام یک دیده می‌شود؟ مارس 2009 و 2011'
$ws.Rows.Item(63).AutoFit()
$ws.Range("A64").Value = 'در گلومرولونفریت وگنر، ویژگی مشخصه مشاهده شده چیست؟'
$ws.Rows.Item(64).AutoFit()
$ws.Range("A65").Value = 'مهم ترین ویژگی سیمان های پلی آکریلیک اسیدی چیست؟'
$ws.Rows.Item(65).AutoFit()
$ws.Range("A66").Value = 'dna میتوکندریایی چیست؟'
$ws.Rows.Item(66).AutoFit()
$ws.Range("A67").Value = 'کدام یک از موارد زیر در مورد ترمیم زخم صحیح است؟'
$ws.Rows.Item(67).AutoFit()
$ws.Range("A68").Value = 'یک زن ۷۰ ساله با چاقی متوسط، با فتق برشی بزرگ در خط وسط مراجعه می‌کند. یک سال قبل تحت عمل رزکسیون کولون به دلیل آدنوکارسینوما قرار گرفته است. کدام یک از گزینه‌های زیر صحیح است؟'
$ws.Rows.Item(68).AutoFit()
$ws.Range("A69").Value = 'مهم‌ترین منبع خونرسانی به معده کدام است؟'
$ws.Rows.Item(69).AutoFit()
$ws.Range("A70").Value = 'موکوپلی ساکاریدی که باقیمانده اسید اورونیک ندارد کدام است؟'
$ws.Rows.Item(70).AutoFit()
$ws.Range("A71").Value = 'توزیع جریان خون عمدتاً توسط چه چیزی تنظیم می‌شود؟'
$ws.Rows.Item(71).AutoFit()
$ws.Range("A72").Value = 'شایع‌ترین علت آمبلیوپی (تنبلی چشم) چیست؟'
$ws.Rows.Item(72).AutoFit()
$ws.Range("A73").Value = 'کدام یک از موارد زیر می‌تواند به ورید اجوف تحتانی (ivc) تهاجم کند؟'
$ws.Rows.Item(73).AutoFit()
$ws.Range("A74").Value = 'سطح تئوفیلین در تمام موارد زیر افزایش می‌یابد، به جز'
$ws.Rows.Item(74).AutoFit()
$ws.Range("A75").Value = 'کدام یک از ساختارهای پریکارد نسبت به درد حساس نیست؟'
$ws.Rows.Item(75).AutoFit()
$ws.Range("A76").Value = 'بیماری پمپه به دلیل کمبود کدام آنزیم در بیماری ذخیره‌ای ایجاد می‌شود؟'
$ws.Rows.Item(76).AutoFit()
$ws.Range("A77").Value = 'علت ایجاد شبه لنگش چیست؟'
$ws.Rows.Item(77).AutoFit()
$ws.Range("A78").Value = 'یک مرد ۲۳ ساله رابطه جنسی محافظت نشده با یک کارگر جنسی تجاری داشت. دو هفته بعد، یک زخم بدون درد و سفت روی سر آلت ایجاد شد که با فشار، سرم شفاف ترشح می‌کرد. غدد لنفاوی کشاله ران در هر دو طرف بزرگ شده و بدون درد بودند. مناسب‌ترین آزمایش تشخیصی چیست؟'
$ws.Rows.Item(78).AutoFit()
$ws.Range("A79").Value = 'دو آنتی‌هیستامین ترفنادین و آستمیزول به دنبال بروز آریتمی‌های قلبی زمانی که در سطح بالایی در خون وجود داشتند، از بازار جمع‌آوری شدند. این اثرات با کدام یک از موارد زیر توضیح داده می‌شود؟'
$ws.Rows.Item(79).AutoFit()
$ws.Range("A80").Value = 'در مورد کلومیفن سیترات کدام گزینه صحیح است؟'
$ws.Rows.Item(80).AutoFit()
$ws.Range("A81").Value = 'ظاهر حباب صابون در انتهای دیستال رادیوس تشخیصی کدام بیماری است؟'
$ws.Rows.Item(81).AutoFit()
$ws.Range("A82").Value = 'کدام یک از عبارات زیر در مورد هاپتن صحیح است؟'
$ws.Rows.Item(82).AutoFit()
$ws.Range("A83").Value = 'درمان رتینوپاتی دیابتی پیشرفته پرولیفراتیو با فیبروز وترئوریتینال گسترده و جداشدگی کششی شبکیه شامل همه موارد زیر به جز کدام است؟'
$ws.Rows.Item(83).AutoFit()
$ws.Range("A84").Value = 'همه موارد زیر در مورد دلیریوم ترمنس صحیح است به جز'
$ws.Rows.Item(84).AutoFit()
$ws.Range("A85").Value = 'dutch cap چیست؟'
$ws.Rows.Item(85).AutoFit()
$ws.Range("A86").Value = 'اگزودای سخت در همه موارد زیر دیده می‌شود به جز –'
$ws.Rows.Item(86).AutoFit()
$ws.Range("A87").Value = 'یک مرد 23 ساله با بیرون زدگی پیشرونده چشم راست مراجعه می‌کند. بیرون زدگی با خم شدن به جلو افزایش می‌یابد و قابل فشرده شدن است. هیچ لرزش یا صدای غیرعادی وجود ندارد. سونوگرافی مداری یک توده اکوژنیک با کانون‌های سایه اکو نشان داد. محتمل‌ترین تشخیص چیست؟'
$ws.Rows.Item(87).AutoFit()
$ws.Range("A88").Value = 'لورازپام به دلیل کدام ویژگی در درمان وضعیت صرعی نسبت به دیازپام ارجحیت دارد؟'
$ws.Rows.Item(88).AutoFit()
$ws.Range("A89").Value = 'کدام یک از عبارات زیر در مورد هایپرپلازی پایدار زجاجیه اولیه (phpv) صحیح نیست؟'
$ws.Rows.Item(89).AutoFit()
$ws.Range("A90").Value = 'شایع‌ترین محل شکستگی کانترکوب در سقوط پس‌سری کدام است؟'
$ws.Rows.Item(90).AutoFit()
$ws.Range("A91").Value = 'پروب پریودنتال سازمان بهداشت جهانی (پروب cpitn) دارای حلقه های رنگی سیاه است'
$ws.Rows.Item(91).AutoFit()
$ws.Range("A92").Value = 'مرز فوقانی سوراخ اپیپلوئیک توسط کدام ساختار تشکیل شده است؟'
$ws.Rows.Item(92).AutoFit()
$ws.Range("A93").Value = 'دوره پنجره به عنوان مدت زمان از ... تعریف می‌شود:'
$ws.Rows.Item(93).AutoFit()
$ws.Range("A94").Value = 'همه موارد زیر در sle دیده می‌شوند به جز -'
$ws.Rows.Item(94).AutoFit()
$ws.Range("A95").Value = 'فریونودرمیا ناشی از کمبود کدام ماده است؟'
$ws.Rows.Item(95).AutoFit()
$ws.Range("A96").Value = 'کدام یک از حالت‌های تهویه زیر بالاترین خطر ایجاد آلکالوز تنفسی را برای بیمار دارد؟'
$ws.Rows.Item(96).AutoFit()
$ws.Range("A97").Value = 'شیگلا را می‌توان از e.coli با تمام ویژگی‌های زیر به جز کدام یک تفکیک کرد؟'
$ws.Rows.Item(97).AutoFit()
$ws.Range("A98").Value = 'بوی نامطبوع و توهم بویایی در آسیب کدام ناحیه دیده می‌شود؟'
$ws.Rows.Item(98).AutoFit()
$ws.Range("A99").Value = 'حرارت مرطوب همه موارد زیر را از بین می برد به جز-'
$ws.Rows.Item(99).AutoFit()
$ws.Range("A100").Value = 'اولین بیحس کننده موضعی که در بیهوشی بالینی استفاده شد کدام بود؟'
$ws.Rows.Item(100).AutoFit()
$ws.Range("A101").Value = 'همه موارد زیر از متغیرهایی هستند که به طور مداوم مرگ و میر ۶ هفته ای را پیش‌بینی می‌کنند، به جز:'
$ws.Rows.Item(101).AutoFit()
$ws.Range("A102").Value = 'مجرای توراسیک را می‌توان با کدام مورد شناسایی کرد؟'
$ws.Rows.Item(102).AutoFit()
$ws.Range("A103").Value = 'ماساژ و استفاده از ضمادها بر روی نواحی دردناک بدن، درد را به دلیل چه چیزی تسکین می‌دهد؟'
$ws.Rows.Item(103).AutoFit()
$ws.Range("A104").Value = 'کدام یک از دندان‌های پیشین زیر اولین دندان‌هایی هستند که در یک نوزاد می‌رویند؟'
$ws.Rows.Item(104).AutoFit()
$ws.Range("A105").Value = 'رأس ریه انسان در حالت ایستاده در مقایسه با پایه آن دارای چه ویژگی است؟'
$ws.Rows.Item(105).AutoFit()
$ws.Range("A106").Value = 'تمام اسیدهای آمینه گلوکوژنیک به جز کدام یک هستند؟'
$ws.Rows.Item(106).AutoFit()
$ws.Range("A107").Value = 'سینه خاموش در کدام مورد مشاهده می‌شود؟'
$ws.Rows.Item(107).AutoFit()
$ws.Range("A108").Value = 'در مراقبت‌های بهداشتی اولیه (phc) چگونه جفت را دفع کنیم؟'
$ws.Rows.Item(108).AutoFit()
$ws.Range("A109").Value = 'یک مرد ۶۰ ساله با مشکل پیشرونده در بلع، استفراغ و برگشت گاه به گاه محتویات معده به مدت ۳ ماه مراجعه کرده است. مطالعات باریم نشان دهنده اتساع قابل توجه مری فوقانی با تنگی در بخش تحتانی است. مانومتری نشان دهنده عدم وجود پریستالسیس مری است. پاتوژنز این وضعیت به احتمال زیاد مربوط به کدام مورد زیر است؟'
$ws.Rows.Item(109).AutoFit()
$ws.Range("A110").Value = 'مجازات تجاوز جنسی تحت کدام بخش از قانون کیفری هند (ipc) قرار دارد؟'
$ws.Rows.Item(110).AutoFit()
$ws.Range("A111").Value = 'همه موارد زیر در مورد کارآزمایی کنترل شده تصادفی صحیح است به جز؟'
$ws.Rows.Item(111).AutoFit()
$ws.Range("A112").Value = 'شایع‌ترین محل شکستگی فک پایین:'
$ws.Rows.Item(112).AutoFit()
$ws.Range("A113").Value = 'گیرنده‌های toll-like مرتبط با ویروس کدامند؟'
$ws.Rows.Item(113).AutoFit()
$ws.Range("A114").Value = 'همه موارد زیر در مورد انعقاد خون صحیح هستند، به جز'
$ws.Rows.Item(114).AutoFit()
$ws.Range("A115").Value = 'کدام یک از موارد زیر از ویژگی‌های نارسایی حاد کلیوی نیست؟'
$ws.Rows.Item(115).AutoFit()
$ws.Range("A116").Value = 'در کدام سیستم، شمارش مداوم تولدها و مرگ‌ها توسط شمارشگر و یک بررسی مستقل توسط سرپرست محقق انجام می‌شود؟'
$ws.Rows.Item(116).AutoFit()
$ws.Range("A117").Value = 'تمام موارد زیر از ویژگی‌های یوئیت قدامی هستند به جز'
$ws.Rows.Item(117).AutoFit()
$ws.Range("A118").Value = 'یک خانم 65 ساله که 7 روز پس از عمل جراحی شکستگی فمور به طور ناگهانی در حال صرف ناهار در بخش فوت می‌کند. محتمل‌ترین علت چیست؟'
$ws.Rows.Item(118).AutoFit()
$ws.Range("A119").Value = 'در مورد آبله مرغان موارد زیر صحیح است به جز:'
$ws.Rows.Item(119).AutoFit()
$ws.Range("A120").Value = 'سندرم ورزشکار با کدام مورد مشخص می‌شود؟'
$ws.Rows.Item(120).AutoFit()
$ws.Range("A121").Value = 'اگر پروتئین‌های سلولی به شکل خاصی تا نشوند، عملکرد آن‌ها تحت تأثیر قرار می‌گیرد. برخی اختلالات در صورت تا نشدن صحیح پروتئین‌های خاص ایجاد می‌شوند. کدام یک از اختلالات زیر به دلیل ایزومریزاسیون ساختاری ایجاد می‌شود؟'
$ws.Rows.Item(121).AutoFit()
$ws.Range("A122").Value = 'کدام عامل بلوک عصبی-عضلانی می‌تواند از جفت عبور کند؟'
$ws.Rows.Item(122).AutoFit()
$ws.Range("A123").Value = 'کوردوما معمولاً کدام نواحی را درگیر می‌کند؟  
الف) ستون فقرات پشتی  
ب) کلیووس  
ج) ستون فقرات کمری  
د) ساکروم  
ه) ستون فقرات گردنی'
$ws.Rows.Item(123).AutoFit()
$ws.Range("A124").Value = 'تیلوزیس به چه چیزی اشاره دارد؟'
$ws.Rows.Item(124).AutoFit()
$ws.Range("A125").Value = 'بلاک عصبی آلوئولار تحتانی به تنهایی می‌تواند در کدام مورد استفاده شود؟'
$ws.Rows.Item(125).AutoFit()
$ws.Range("A126").Value = 'کودکی با هایپرآمونمی، بعدها دچار پانکراتیت و سکته گانگلیون پایه می‌شود. ممکن است مبتلا به کدام یک از موارد زیر باشد؟'
$ws.Rows.Item(126).AutoFit()
$ws.Range("A127").Value = 'همه موارد زیر در مورد اسیدوز توبولی کلیوی صحیح هستند به جز.'
$ws.Rows.Item(127).AutoFit()
$ws.Range("A128").Value = 'cephalhematoma:'
$ws.Rows.Item(128).AutoFit()
$ws.Range("A129").Value = 'یک مرد 65 ساله سیگاری با هماچوری کامل بدون درد مراجعه می‌کند. محتمل‌ترین تشخیص چیست؟'
$ws.Rows.Item(129).AutoFit()
$ws.Range("A130").Value = 'عصب مدیان همه موارد زیر را عصب‌دهی می‌کند به جز:'
$ws.Rows.Item(130).AutoFit()
$ws.Range("A131").Value = 'یک بیمار 45 ساله با کاهش شنوایی هدایتی، درد صورت در ناحیه گیجگاهی-آهیانه‌ای و فک پایین، و عدم تحرک کام نرم در سمت راست مراجعه می‌کند. تشخیص احتمالی چیست؟'
$ws.Rows.Item(131).AutoFit()
$ws.Range("A132").Value = 'شایع‌ترین محل باز شدن داخلی فیستول برونشی در کجاست؟'
$ws.Rows.Item(132).AutoFit()
$ws.Range("A133").Value = 'کدام یک از روغن های خوراکی زیر بیشترین مقدار اسیدهای چرب چند غیراشباع را تولید می کند؟'
$ws.Rows.Item(133).AutoFit()
$ws.Range("A134").Value = 'علاوه بر 3 d (اسهال، درماتیت و زوال عقل) در کمبود نیاسین، d چهارم نشان دهنده چیست؟'
$ws.Rows.Item(134).AutoFit()
$ws.Range("A135").Value = 'قانون نیستن مربوط به کدام مورد است؟'
$ws.Rows.Item(135).AutoFit()
$ws.Range("A136").Value = 'کدام سوبسترا در دوره پس از عمل زودتر تخلیه می‌شود؟'
$ws.Rows.Item(136).AutoFit()
$ws.Range("A137").Value = 'شایع‌ترین عارضه سرخک کدام است؟'
$ws.Rows.Item(137).AutoFit()
$ws.Range("A138").Value = 'توده در حفره ایلیاک راست می‌تواند ناشی از چه باشد؟'
$ws.Rows.Item(138).AutoFit()
$ws.Range("A139").Value = 'اگر ضایعه در کدام ناحیه باشد، مکولار اسپارینگ دیده می‌شود؟'
$ws.Rows.Item(139).AutoFit()
$ws.Range("A140").Value = 'روش مطالعه مورد-شاهدی دارای سه ویژگی متمایز زیر است. گزینه نامناسب را انتخاب کنید'
$ws.Rows.Item(140).AutoFit()
$ws.Range("A141").Value = 'بواسیزوماب در چه موردی استفاده می‌شود؟'
$ws.Rows.Item(141).AutoFit()
$ws.Range("A142").Value = 'یک مرد ۴۷ ساله با سابقه بیماری سلول داسی شکل، چندین بار بستری شده و نیاز به قرار دادن خطوط داخل وریدی داشته است. بیمار دسترسی وریدی محیطی ضعیفی دارد و یک کاتتر در ورید سابکلاوین راست قرار داده می‌شود. بیمار متعاقباً دچار ناراحتی و تورم در بازوی راست و تب ۴۰.۱ درجه سانتی‌گراد با لرز می‌شود. چندین کشت خون گرفته می‌شود و کوکسی‌های گرم مثبت جدا می‌شوند. این ارگانیسم کاتالاز مثبت است و روی آگار مانیتول نمک رشد می‌کند، اما آگار را زرد نمی‌کند؛ کلونی‌ها روی پلیت آگار خون گوسفند گاما-همولیتیک هستند. کدام یک از ارگانیسم‌های زیر به احتمال زیاد باعث علائم این بیمار شده است؟'
$ws.Rows.Item(142).AutoFit()
$ws.Range("A143").Value = 'سیاهی رفتن چشم به دلیل کدام مورد است؟'
$ws.Rows.Item(143).AutoFit()
$ws.Range("A144").Value = 'گرسنگی و دیابت شیرین می‌توانند منجر به کتوزیس شوند. کدام یک از ویژگی‌های زیر در کتوزیس ناشی از دیابت شیرین وجود دارد یا مربوط به آن است؟'
$ws.Rows.Item(144).AutoFit()
$ws.Range("A145").Value = 'کدام یک از موارد زیر نشانه قطعی نارسایی قلبی در نوزاد مبتلا به بیماری مادرزادی قلب است؟'
$ws.Rows.Item(145).AutoFit()
$ws.Range("A146").Value = 'کم خونی اسپور سل در کدام مورد مشاهده می‌شود؟'
$ws.Rows.Item(146).AutoFit()
$ws.Range("A147").Value = 'خونریزی شدید از مقعد در یک بیمار 70 ساله به چه علتی است؟'
$ws.Rows.Item(147).AutoFit()
$ws.Range("A148").Value = 'کدام یک از ویروس‌های زیر بدون پوشش است اما حاوی rna دو رشته‌ای می‌باشد؟'
$ws.Rows.Item(148).AutoFit()
$ws.Range("A149").Value = 'شاخه کلیوی ورید اجوف تحتانی از کدام یک تشکیل می‌شود؟'
$ws.Rows.Item(149).AutoFit()
$ws.Range("A150").Value = 'شایع‌ترین عارضه جراحی پس از عمل وایپل چیست؟'
$ws.Rows.Item(150).AutoFit()
$ws.Range("A151").Value = 'در مقایسه با یک جمعیت ناهمگن، میزان بروز مال اکلوژن در یک جمعیت همگن عموماً:'
$ws.Rows.Item(151).AutoFit()
$ws.Range("A152").Value = 'سوءتغذیه مزمن با کدام شاخص اندازه‌گیری می‌شود؟'
$ws.Rows.Item(152).AutoFit()
$ws.Range("A153").Value = 'این سونوگرافی اندازه‌گیری crl جنین ۸ هفته‌ای را نشان می‌دهد. پیکان در تصویر به کدام یک از موارد زیر اشاره دارد؟'
$ws.Rows.Item(153).AutoFit()
$ws.Range("A154").Value = 'همه موارد زیر در مورد فلج دوطرفه دیافراگم صحیح هستند، به جز:'
$ws.Rows.Item(154).AutoFit()
$ws.Range("A155").Value = 'پروتئینوری توبول پروگزیمال و ضایعات استخوانی دردناک در مصرف بیش از حد کدام ماده دیده می‌شود؟'
$ws.Rows.Item(155).AutoFit()
$ws.Range("A156").Value = 'در فاز تاریکی چرخه بینایی، کدام فرم ویتامین a با اپسین ترکیب می‌شود تا رودوپسین را بسازد؟'
$ws.Rows.Item(156).AutoFit()
$ws.Range("A157").Value = 'رادیوگرافی زیر نشان‌دهنده کدام مورد است؟'
$ws.Rows.Item(157).AutoFit()
$ws.Range("A158").Value = 'لکه‌های فاکس فوردایس نشانگر چه چیزی هستند؟'
$ws.Rows.Item(158).AutoFit()
$ws.Range("A159").Value = 'برای آدنوکارسینوم انتهای تحتانی مری، ازوفاژکتومی ترانسهیاتال برنامهریزی شده است. رویکرد به ترتیب زیر خواهد بود:'
$ws.Rows.Item(159).AutoFit()
$ws.Range("A160").Value = 'شایع‌ترین محل کیست تیروگلوسال کجاست؟'
$ws.Rows.Item(160).AutoFit()
$ws.Range("A161").Value = 'تست کالریک برای چه چیزی انجام می‌شود؟'
$ws.Rows.Item(161).AutoFit()
$ws.Range("A162").Value = 'سوراخ شدن خلفی زخم معده به کجا تخلیه می‌شود؟'
$ws.Rows.Item(162).AutoFit()
$ws.Range("A163").Value = 'کدام یک از موارد زیر از پرولاپس رحم جلوگیری نمی‌کند؟'
$ws.Rows.Item(163).AutoFit()
$ws.Range("A164").Value = 'سموم سینرژی‌هیمنوتروپیک استافیلوکوکی شامل کدام موارد است؟
الف) سم آلفا
ب) سم بتا
ج) سم گاما
د) سم دلتا
ه) سم پانتون-والنتاین'
$ws.Rows.Item(164).AutoFit()
$ws.Range("A165").Value = 'تمام موارد زیر در مورد هموسیدروز ریوی ایدیوپاتیک صحیح هستند به جز:'
$ws.Rows.Item(165).AutoFit()
$ws.Range("A166").Value = 'کدام عبارت این جمله را به بهترین شکل تکمیل می‌کند؟ کانال اینگوینال:'
$ws.Rows.Item(166).AutoFit()
$ws.Range("A167").Value = 'یک بیمار از منطقه قبیله‌ای جارکند با تب به مدت ۳ روز مراجعه کرده است. نمونه خون محیطی جمع‌آوری و با گیمسا رنگ‌آمیزی شده است. تشخیص مالاریا داده شده است. اسمیر در شکل نشان داده شده است. علت احتمالی چیست؟'
$ws.Rows.Item(167).AutoFit()
$ws.Range("A168").Value = 'یکی از شکستگی‌های شایعی که در حین بوکس به دلیل ضربه با مشت بسته رخ می‌دهد، کدام است؟'
$ws.Rows.Item(168).AutoFit()
$ws.Range("A169").Value = 'کدام یک از انواع سلول‌های زیر به طور خاص مربوط به عفونت نهفته تناسلی با hsv-2 هستند؟ (غیر مرتبط)'
$ws.Rows.Item(169).AutoFit()
$ws.Range("A170").Value = 'تمام موارد زیر کیستهای تکا لوتئین را از کیستهای جسم زرد متمایز میکنند، به جز'
$ws.Rows.Item(170).AutoFit()
$ws.Range("A171").Value = 'کرایوپرسیپیتات در کدام مورد مفید است؟'
$ws.Rows.Item(171).AutoFit()
$ws.Range("A172").Value = 'مثلث کالو یک نشانگر مهم برای کدام ساختار است؟'
$ws.Rows.Item(172).AutoFit()
$ws.Range("A173").Value = 'پارامترهای ایده‌آل برای ماساژ قلبی در احیای قلبی ریوی شامل همه موارد زیر به جز کدام است؟'
$ws.Rows.Item(173).AutoFit()
$ws.Range("A174").Value = 'میزان نیاز آهن و اسید فولیک برای کودکان 6 ماهه تا 10 ساله چقدر است؟'
$ws.Rows.Item(174).AutoFit()
$ws.Range("A175").Value = 'سیروز میکرونودولار در همه موارد زیر دیده می‌شود به جز'
$ws.Rows.Item(175).AutoFit()
$ws.Range("A176").Value = 'برش‌های آزمایشی در کدام یک از موارد زیر دیده می‌شوند؟ (pgi june 2009)'
$ws.Rows.Item(176).AutoFit()
$ws.Range("A177").Value = 'ظاهر شدن «رگه‌های مرمری» در یک جسد چقدر طول می‌کشد؟'
$ws.Rows.Item(177).AutoFit()
$ws.Range("A178").Value = 'کدام یک از سمان های زیر بیشترین تحریک را برای بافت پالپ دارد؟'
$ws.Rows.Item(178).AutoFit()
$ws.Range("A179").Value = 'استئوسارکوم ثانویه با کدام یک از موارد زیر مرتبط است؟'
$ws.Rows.Item(179).AutoFit()
$ws.Range("A180").Value = 'پالپ کپ چیست؟'
$ws.Rows.Item(180).AutoFit()
$ws.Range("A181").Value = 'سلول‌های ادونتوبلاست جدید از سلول‌های مزانشیمی در چه زمانی تمایز می‌یابند؟'
$ws.Rows.Item(181).AutoFit()
$ws.Range("A182").Value = 'یک مرد 49 ساله مبتلا به سرطان پروستات تحت اشعه ایکس قرار گرفت. در تصاویر اشعه ایکس، مناطق اسکلروز و فروپاشی مهره‌های t10 و t11 مشاهده شد. گسترش این سرطان به مهره‌های فوق از طریق کدام مسیر صورت گرفته است؟'
$ws.Rows.Item(182).AutoFit()
$ws.Range("A183").Value = 'علامت t در کدام مورد دیده می‌شود؟'
$ws.Rows.Item(183).AutoFit()
$ws.Range("A184").Value = 'آپاندیس بیضه باقیمانده کدام است؟'
$ws.Rows.Item(184).AutoFit()
$ws.Range("A185").Value = 'یک بیمار مرد 45 ساله با شکایت اصلی فراموشی پیشرونده در طی چند سال گذشته همراه با آتاکسی به درمانگاه پزشکی مراجعه کرد. در معاینه: راه رفتن با پایه عریض و شافلینگ، سفتی عمومی، هایپراکستنشن گردن. در معاینه چشم: محدودیت حرکات چشمی، ساکادهای چشمی آهسته، اختلال در نگاه به پایین. mri مغز انجام شد. در این شرایط کدام پروتئین بیشتر در ساختارهای سابکورتیکال رسوب می‌کند؟'
$ws.Rows.Item(185).AutoFit()
$ws.Range("A186").Value = 'گال، یک عفونت پوستی ناشی از sarcoptes scabiei، نمونه ای از کدام مورد است؟'
$ws.Rows.Item(186).AutoFit()
$ws.Range("A187").Value = 'همه موارد زیر در مورد آرتریت تمپورال صحیح است به جز -'
$ws.Rows.Item(187).AutoFit()
$ws.Range("A188").Value = 'عملکرد عضلات بین دنده ای خارجی چیست؟'
$ws.Rows.Item(188).AutoFit()
$ws.Range("A189").Value = 'کدام یک از موارد زیر توسط سلول های پانت ترشح می شود؟'
$ws.Rows.Item(189).AutoFit()
$ws.Range("A190").Value = 'کودکی با عفونت دستگاه تنفسی، ظرف ۲ روز پس از عفونت دچار هماچوری شدید می‌شود. تشخیص احتمالی چیست؟'
$ws.Rows.Item(190).AutoFit()
$ws.Range("A191").Value = 'پیشگیری برای پرسنل بهداشتی که در بخش طاعون کار می‌کنند چیست؟  
الف) واکسن  
ب) تتراسایکلین در طول مدت خدمت  
ج) یک دوره تتراسایکلین  
د) واکسن و اریترومایسین  
ه) مشاهده'
$ws.Rows.Item(191).AutoFit()
$ws.Range("A192").Value = 'کدام یک از داروهای زیر در درمان بیماران مبتلا به ایدز باعث سرکوب مغز استخوان نمی‌شود؟'
$ws.Rows.Item(192).AutoFit()
$ws.Range("A193").Value = 'کدام یک از موارد زیر یک عارضه جانبی جدی ناشی از زولندرونات است؟'
$ws.Rows.Item(193).AutoFit()
$ws.Range("A194").Value = 'علت کدورت مدیاستین خلفی در نمای pa و لترال رادیوگرافی قفسه سینه'
$ws.Rows.Item(194).AutoFit()
$ws.Range("A195").Value = 'تمام مکانیسم‌های زیر ممکن است در کاهش خطر عفونت دستگاه تناسلی فوقانی در کاربران آی‌یودی‌های آزادکننده پروژستین نقش داشته باشند، به جز:'
$ws.Rows.Item(195).AutoFit()
$ws.Range("A196").Value = 'کدام یک از داروهای زیر در درمان هیپرکلسمی با علت ناشناخته استفاده نمی‌شود؟'
$ws.Rows.Item(196).AutoFit()
$ws.Range("A197").Value = 'کدام یک از داروهای ضد سرطان زیر می‌توانند از سد خونی مغزی عبور کنند؟'
$ws.Rows.Item(197).AutoFit()
$ws.Range("A198").Value = 'کدام ویتامین به عنوان کوفاکتور در متابولیسم گلیسین مرتبط است؟'
$ws.Rows.Item(198).AutoFit()
$ws.Range("A199").Value = 'یک زن 38 ساله با اختلال دوقطبی که در طول 2 سال گذشته با لیتیوم پایدار بوده است، با سابقه 2 هفته ای خلق افسرده، تمرکز ضعیف، از دست دادن اشتها و افکار خودکشی غیرفعال به مطب روانپزشک خود مراجعه می کند. روانپزشک باید کدام یک از مراحل زیر را انجام دهد؟'
$ws.Rows.Item(199).AutoFit()
$ws.Range("A200").Value = 'ترومبوآنژئیت انسدادی (بیماری بورگر) کدام موارد را درگیر می‌کند؟'
$ws.Rows.Item(200).AutoFit()
$ws.Range("A201").Value = 'کدام نوع سوراخ شدن پرده صماخ بیشتر در csom توبوتیمپانیک دیده می‌شود؟'
$ws.Rows.Item(201).AutoFit()
$ws.Range("A202").Value = 'آی یو دی مسی 200 به دلیل چه چیزی به عنوان مس 200 شناخته می‌شود؟'
$ws.Rows.Item(202).AutoFit()
$ws.Range("A203").Value = 'داروی انتخابی برای درمان اولیه پارکینسونیسم کدام است؟'
$ws.Rows.Item(203).AutoFit()
$ws.Range("A204").Value = 'همولیز در کمبود آنزیم g6pd (گلوکز 6 فسفات دهیدروژناز) ممکن است با تمام داروهای زیر رخ دهد به جز:'
$ws.Rows.Item(204).AutoFit()
$ws.Range("A205").Value = 'شایع‌ترین علت بزرگ‌شدگی سوپراسلار با کلسیفیکاسیون در کودکان چیست؟'
$ws.Rows.Item(205).AutoFit()
$ws.Range("A206").Value = 'کلسيتونين توسط کدام غده ترشح می‌شود؟'
$ws.Rows.Item(206).AutoFit()
$ws.Range("A207").Value = 'در هنگام تطابق، کاهش در شعاع انحنای کدام قسمت رخ می‌دهد؟'
$ws.Rows.Item(207).AutoFit()
$ws.Range("A208").Value = 'کدام یک از موارد زیر می‌تواند علت یافته سونوگرافی داده شده باشد؟'
$ws.Rows.Item(208).AutoFit()
$ws.Range("A209").Value = 'شایع‌ترین علت روبئوزیس ایریدیس چیست؟'
$ws.Rows.Item(209).AutoFit()
$ws.Range("A210").Value = 'در حالت ایستاده، جریان خون ضربان‌دار در کدام ناحیه از ریه سالم یافت می‌شود؟'
$ws.Rows.Item(210).AutoFit()
$ws.Range("A211").Value = 'نوع تهاجمی لنفوم سلول t پوستی کدام است؟'
$ws.Rows.Item(211).AutoFit()
$ws.Range("A212").Value = 'کدام دارو برای کنترل تاکیکاردی فوق بطنی به عنوان داروی انتخابی در نظر گرفته می‌شود؟'
$ws.Rows.Item(212).AutoFit()
$ws.Range("A213").Value = 'یک پسر ۲.۵ هفته‌ای با استفراغ پرتابی غیرصفراوی به بخش اورژانس کودکان آورده می‌شود. در معاینه فیزیکی، یک توده سفت و زیتونی‌شکل در ناحیه اپیگاستر مشاهده می‌شود و یافته‌های سونوگرافی موارد زیر را نشان می‌دهد. کودک تحت درمان جراحی قرار می‌گیرد. یک هفته بعد، کودک مجدداً با اسهال شل پس از تغذیه با شیر مادر مراجعه می‌کند. علت احتمالی این دوره اسهال چیست؟'
$ws.Rows.Item(213).AutoFit()
$ws.Range("A214").Value = 'کدام یک از سیستم‌های حسی زیر از فیبرهای بدون میلین برای انتقال اطلاعات به سیستم عصبی مرکزی استفاده می‌کند؟'
$ws.Rows.Item(214).AutoFit()
$ws.Range("A215").Value = 'لایزین در کدام یک کمبود دارد؟'
$ws.Rows.Item(215).AutoFit()
$ws.Range("A216").Value = 'یک نجار 39 ساله دو بطری مشروب از مغازه محلی خریداری کرده است. پس از حدود یک ساعت، دچار گیجی، استفراغ و تاری دید شده است. او به بخش اورژانس آورده شده است. باید به او داده شود:'
$ws.Rows.Item(216).AutoFit()
$ws.Range("A217").Value = 'کدام یک از داروهای زیر در مدیریت خونریزی پس از زایمان استفاده نمی‌شود؟'
$ws.Rows.Item(217).AutoFit()
$ws.Range("A218").Value = 'تعداد کروموزوم‌ها در سندرم ترنر -'
$ws.Rows.Item(218).AutoFit()
$ws.Range("A219").Value = 'کدام یک از موارد زیر نشان‌دهنده شکست در انتقال بیمار از تهویه کمکی به تهویه داوطلبانه است؟'
$ws.Rows.Item(219).AutoFit()
$ws.Range("A220").Value = 'کدام یک از موارد زیر باعث هیدرولیز پپتیدوگلیکان ها می شود؟'
$ws.Rows.Item(220).AutoFit()
$ws.Range("A221").Value = 'شکستگی هر دو استخوان ساعد در یک سطح، وضعیت بازو در گچ به چه صورت است؟'
$ws.Rows.Item(221).AutoFit()
$ws.Range("A222").Value = 'ظاهر کاغذ مچاله شده سلول‌ها ویژگی کدام بیماری است؟'
$ws.Rows.Item(222).AutoFit()
$ws.Range("A223").Value = 'عبارت صحیح در مورد آسیب حاد کلیوی ناشی از کنتراست کدام است؟'
$ws.Rows.Item(223).AutoFit()
$ws.Range("A224").Value = 'روش تشخیصی انتخابی برای سرطان پستان پرخطر در زنان کدام است؟'
$ws.Rows.Item(224).AutoFit()
$ws.Range("A225").Value = 'شایع ترین علت آترزی کوانال چیست؟'
$ws.Rows.Item(225).AutoFit()
$ws.Range("A226").Value = 'عبارت نادرست در مورد فلج مغزی کدام است؟'
$ws.Rows.Item(226).AutoFit()
$ws.Range("A227").Value = 'تومور سلول زایای دوطرفه چیست؟'
$ws.Rows.Item(227).AutoFit()
$ws.Range("A228").Value = 'یک بیمار در بخش مراقبت‌های ویژه فشار طبیعی شریان ریوی دارد؛ مقاومت عروق محیطی سیستمیک پایین است؛ شاخص قلبی پایین است؛ فشار اکسیژن شریانی 93 است. احتمالات تشخیصی کدام‌اند؟'
$ws.Rows.Item(228).AutoFit()
$ws.Range("A229").Value = 'بیماری جنگل کیسانور توسط چه چیزی ایجاد می‌شود؟'
$ws.Rows.Item(229).AutoFit()
$ws.Range("A230").Value = 'سیاه‌شدن حافظه ناشی از الکل در چه میزان از مسمومیت با الکل مشاهده می‌شود؟'
$ws.Rows.Item(230).AutoFit()
$ws.Range("A231").Value = 'در مورد بنزودیازپین‌ها کدام گزینه صحیح است؟'
$ws.Rows.Item(231).AutoFit()
$ws.Range("A232").Value = 'در مورد رتینوبلاستوما همه موارد زیر صحیح هستند به جز:'
$ws.Rows.Item(232).AutoFit()
$ws.Range("A233").Value = 'آینه حنجره قبل از استفاده با قرار دادن چه چیزی گرم می‌شود؟'
$ws.Rows.Item(233).AutoFit()
$ws.Range("A234").Value = 'کدام یک از داروهای زیر عمدتاً برای سیستم دهلیزی سمی نیست؟'
$ws.Rows.Item(234).AutoFit()
$ws.Range("A235").Value = 'همه باعث افزایش وزن می‌شوند، به جز:'
$ws.Rows.Item(235).AutoFit()
$ws.Range("A236").Value = 'نوع شکستگی مونتگیا که معمولاً با آسیب های عصبی همراه است کدام است؟'
$ws.Rows.Item(236).AutoFit()
$ws.Range("A237").Value = 'کدام یک جزو هسته‌های مخچه نیست؟'
$ws.Rows.Item(237).AutoFit()
$ws.Range("A238").Value = 'سوماتواستاتین توسط کدام سلول ترشح می‌شود؟'
$ws.Rows.Item(238).AutoFit()
$ws.Range("A239").Value = 'در لارنژکتومی جزئی حنجره، تمامی ساختارهای زیر به جز کدام یک برداشته می‌شوند؟'
$ws.Rows.Item(239).AutoFit()
$ws.Range("A240").Value = 'اکتینومایکوزیس به کدام یک از موارد زیر حساس است؟'
$ws.Rows.Item(240).AutoFit()
$ws.Range("A241").Value = 'در مردی که چمدانی را بلند می‌کند، دررفتگی خلفی مفصل گلنوهومرال توسط کدام عضله جلوگیری می‌شود؟'
$ws.Rows.Item(241).AutoFit()
$ws.Range("A242").Value = 'لینیتیس پلاستیکا چیست؟'
$ws.Rows.Item(242).AutoFit()
$ws.Range("A243").Value = 'همه موارد زیر از ویژگی‌های شبه‌تومور مغزی هستند به جز:'
$ws.Rows.Item(243).AutoFit()
$ws.Range("A244").Value = 'شایع‌ترین بیماری ارثی نابینایی ناشی از ناهنجاری کروموزومی میتوکندریایی کدام است؟'
$ws.Rows.Item(244).AutoFit()
$ws.Range("A245").Value = 'اثر کورتیزول'
$ws.Rows.Item(245).AutoFit()
$ws.Range("A246").Value = 'پشت بدن کودکان 10 تا 14 ساله چه درصدی از کل سطح بدن را تشکیل می‌دهد؟'
$ws.Rows.Item(246).AutoFit()
$ws.Range("A247").Value = 'یک دختر نوزده ساله با قد کوتاه، نوک پستان‌های با فاصله زیاد و آمنوره اولیه به احتمال زیاد دارای کاریوتایپ زیر است -'
$ws.Rows.Item(247).AutoFit()
$ws.Range("A248").Value = 'کدام یک از موارد زیر در عکس رادیوگرافی قفسه سینه به طور بصری تشخیصی برای آئوریتیس است؟'
$ws.Rows.Item(248).AutoFit()
$ws.Range("A249").Value = 'cyclic gmp بر روی کدام یک اثر می‌کند؟'
$ws.Rows.Item(249).AutoFit()
$ws.Range("A250").Value = 'نمره سیلورمن را برای یک نوزاد با تنفس پارادوکسیکال، عقب نشینی خفیف قفسه سینه پایین، عقب نشینی واضح زائده گزیفوئید، باز شدن کم بینی و خرخر بازدمی قابل شنیدن بدون استتوسکوپ محاسبه کنید'
$ws.Rows.Item(250).AutoFit()
$ws.Range("A251").Value = 'اصطلاح "id" در مدل ساختاری ذهن به چه چیزی اشاره دارد؟'
$ws.Rows.Item(251).AutoFit()
$ws.Range("A252").Value = 'کدام یک از موارد زیر در هنگام لوله گذاری داخل تراشه منع مصرف دارد؟ مارس 2011'
$ws.Rows.Item(252).AutoFit()
$ws.Range("A253").Value = 'نوزاد هیچ تلاشی برای تنفس ندارد، هیچ پاسخی به تحریک یا رفلکس نشان نمی‌دهد، ضربان قلب کمتر از 100 است و اندام‌های انتهایی آبی رنگ هستند. نمره کلی آپگار ________ است.'
$ws.Rows.Item(253).AutoFit()
$ws.Range("A254").Value = 'کدام نوع آدنوم بیشترین تمایل را به کلسیفیکاسیون دیستروفیک دارد؟'
$ws.Rows.Item(254).AutoFit()
$ws.Range("A255").Value = 'چگونه می‌توان پاسخ به درمان ویتامین d در راشیتیسم را ارزیابی کرد؟'
$ws.Rows.Item(255).AutoFit()
$ws.Range("A256").Value = 'تمامی موارد زیر در میتوز رخ می‌دهند به جز'
$ws.Rows.Item(256).AutoFit()
$ws.Range("A257").Value = 'همه موارد زیر با ضعف "اپیزودیک" مرتبط هستند، به جز -'
$ws.Rows.Item(257).AutoFit()
$ws.Range("A258").Value = 'کلیه تولید می‌کند'
$ws.Rows.Item(258).AutoFit()
$ws.Range("A259").Value = 'آسپرین با کدام مکانیسم باعث آسیب به مخاط می‌شود؟'
$ws.Rows.Item(259).AutoFit()
$ws.Range("A260").Value = 'شایع ترین عامل ایجاد کننده مننژیت نوزادی کدام است؟'
$ws.Rows.Item(260).AutoFit()
$ws.Range("A261").Value = 'طولانی‌ترین دوره تأخیر در رشد باکتریایی با کدام یک از موارد زیر مرتبط است؟'
$ws.Rows.Item(261).AutoFit()
$ws.Range("A262").Value = 'یک زن ۲۰ ساله مبتلا به بیماری هیرشپرونگ با درد حاد پا مراجعه می‌کند. بیمار ۳ سال پیش یک گلیوما برداشته شده است. تصویربرداری اشعه ایکس از پا نشان‌دهنده شکستگی درشت‌نی چپ است. آزمایش‌های آزمایشگاهی سطح بالای کلسیم و pth در سرم را نشان می‌دهند. سی‌تی‌اسکن گردن یک توده منفرد پاراتیروئید را نشان می‌دهد. دو سال بعد، بیمار با فشار خون بالا مراجعه می‌کند و سی‌تی‌اسکن شکم یک توده ۴ سانتی‌متری در غده فوق‌کلیوی راست را نشان می‌دهد. مطالعات ژنتیکی انجام شده روی این بیمار به احتمال زیاد جهش‌های ژرم‌لاین در کدام یک از پروتئوآنکوژن‌های زیر را نشان می‌دهد؟'
$ws.Rows.Item(262).AutoFit()
$ws.Range("A263").Value = 'سطح طبیعی بیلی روبین سرم چقدر است؟'
$ws.Rows.Item(263).AutoFit()
$ws.Range("A264").Value = 'همیانوپیای دو گیجگاهی در ضایعه کدام ناحیه دیده می‌شود؟'
$ws.Rows.Item(264).AutoFit()
$ws.Range("A265").Value = 'کدام یک از عبارات زیر تفاوت بین فکر وسواسی و هذیان را مشخص می کند؟'
$ws.Rows.Item(265).AutoFit()
$ws.Range("A266").Value = 'کدام یک از گزینه های زیر بیشترین ارتباط را با لومفانترین دارد؟'
$ws.Rows.Item(266).AutoFit()
$ws.Range("A267").Value = 'سرطان سلول سنگفرشی مثانه مستعد کدام یک از موارد زیر است؟  
الف) سنگ ادراری  
ب) اوراکوس پایدار  
ج) شیستوزومیازیس  
د) پولیپ  
ه) سیگار کشیدن'
$ws.Rows.Item(267).AutoFit()
$ws.Range("A268").Value = 'پروکتیت رادیاسیون مزمن با درمان تمام بدخیمی‌ها مرتبط است، به جز:'
$ws.Rows.Item(268).AutoFit()
$ws.Range("A269").Value = 'همه موارد زیر می‌توانند آنژین را تشدید کنند به جز'
$ws.Rows.Item(269).AutoFit()
$ws.Range("A270").Value = 'بیماری نوزادان که قابل پیشگیری با واکسن است کدام است؟'
$ws.Rows.Item(270).AutoFit()
$ws.Range("A271").Value = 'یک رزیدنت بیهوشی در حال انجام بیحسی نخاعی بود که بیمار ناگهان دچار بیصدایی و از دست دادن هوشیاری شد. چه اتفاقی ممکن است افتاده باشد؟'
$ws.Rows.Item(271).AutoFit()
$ws.Range("A272").Value = 'کدام یک از گزینه‌های زیر جزو شاخه‌های قوس آئورت نیست؟'
$ws.Rows.Item(272).AutoFit()
$ws.Range("A273").Value = 'نوزادی با بزرگ‌شدگی کلیتوریس و درجاتی از چسبندگی لبیا مشاهده می‌شود. اولین قدم در رسیدن به تشخیص چیست؟'
$ws.Rows.Item(273).AutoFit()
$ws.Range("A274").Value = 'بیماری فون هیپل لینداو با همه موارد زیر مرتبط است به جز،'
$ws.Rows.Item(274).AutoFit()
$ws.Range("A275").Value = 'مسیر پارووسلولار، از هسته زانویی جانبی به قشر بینایی، سیگنال‌های مربوط به تشخیص کدام یک از موارد زیر را حمل می‌کند؟'
$ws.Rows.Item(275).AutoFit()
$ws.Range("A276").Value = 'همه موارد زیر باعث پوکی استخوان می‌شوند به جز'
$ws.Rows.Item(276).AutoFit()
$ws.Range("A277").Value = 'شایع‌ترین علت اوتیت مدیای حاد در کودکان چیست؟'
$ws.Rows.Item(277).AutoFit()
$ws.Range("A278").Value = 'یک بیمار مرد ۵۵ ساله از رشد غیرطبیعی در گونه راست شکایت دارد. بیمار عادت به جویدن تنباکو به مدت ۳۰ سال دارد. یافته‌های رادیولوژیک یک توده قارچ‌مانند ۸x8 سانتی‌متری را نشان می‌دهد که ماگزیلا و مندیبل را در سمت راست مورد تهاجم قرار داده است. به بیمار شیمی‌درمانی تسکینی با سیس‌پلاتین توصیه شد. تمامی تغییرات زیر در dna توسط این دارو ایجاد می‌شوند به جز-'
$ws.Rows.Item(278).AutoFit()
$ws.Range("A279").Value = 'بیضه‌ها به طور کامل در کیسه بیضه تا چه زمانی پایین می‌آیند؟ dnb 09'
$ws.Rows.Item(279).AutoFit()
$ws.Range("A280").Value = 'در رفلکس کوشینگ، کدام مورد مشاهده می‌شود؟'
$ws.Rows.Item(280).AutoFit()
$ws.Range("A281").Value = 'شایع ترین تومور حساس به پرتوی تخمدان کدام است؟'
$ws.Rows.Item(281).AutoFit()
$ws.Range("A282").Value = 'یک زن ۱۹ ساله از نظر جنسی فعال با درد شکم و حالت تهوع به مرکز اورژانس مراجعه می‌کند. آخرین قاعدگی او هفت روز پیش شروع شده است. معاینه فیزیکی تب و درد زیر شکم را نشان می‌دهد. معاینه لگن ترشحات موکوپورولنت، حساسیت آدنکس و حساسیت با حرکت دهانه رحم را نشان می‌دهد. تست بارداری منفی است؛ تعداد گلبول‌های سفید بیمار افزایش یافته است. کدام یک از عوارض زیر با این اختلال مرتبط است؟'
$ws.Rows.Item(282).AutoFit()
$ws.Range("A283").Value = 'شایع‌ترین عامل زمینه‌ساز اریتروپلاکی و لوکوپلاکی کدام است؟'
$ws.Rows.Item(283).AutoFit()
$ws.Range("A284").Value = 'کدام یک از موارد زیر در rasopathies دیده نمی‌شود؟'
$ws.Rows.Item(284).AutoFit()
$ws.Range("A285").Value = 'نمای ترجیحی اشعه ایکس برای کف حفره چشم کدام است؟'
$ws.Rows.Item(285).AutoFit()
$ws.Range("A286").Value = 'شبه زوال عقل در کدام یک از موارد زیر مشاهده می‌شود؟'
$ws.Rows.Item(286).AutoFit()
$ws.Range("A287").Value = 'گنادوتروپین چیست؟'
$ws.Rows.Item(287).AutoFit()
$ws.Range("A288").Value = 'توکسین خارج سلولی کد شده توسط فاژ در ویبریو کلرا شبیه به کدام توکسین اشرشیاکلی است؟'
$ws.Rows.Item(288).AutoFit()
$ws.Range("A289").Value = 'کدام یک از زمان‌های زیر برای ترمیم شکاف کام ایده‌آل است؟'
$ws.Rows.Item(289).AutoFit()
$ws.Range("A290").Value = 'بهترین موقعیت برای فلپ جابجا شده به سمت اپیکال کدام است؟'
$ws.Rows.Item(290).AutoFit()
$ws.Range("A291").Value = 'یک بیمار اسکیزوفرنی تحت درمان با نورولپتیک‌ها، علائم روان‌پریشی او بهبود یافته اما دچار غمگینی شده، کمتر با دیگران صحبت می‌کند و بیشتر در تخت می‌ماند. همه موارد زیر می‌توانند علت این وضعیت باشند به جز:'
$ws.Rows.Item(291).AutoFit()
$ws.Range("A292").Value = 'اتیکت پزشکی چیست؟'
$ws.Rows.Item(292).AutoFit()
$ws.Range("A293").Value = 'تمامی عبارات زیر در مورد واکسن dpt صحیح هستند به جز:'
$ws.Rows.Item(293).AutoFit()
$ws.Range("A294").Value = 'دنده‌ها از چه چیزی تشکیل می‌شوند؟'
$ws.Rows.Item(294).AutoFit()
$ws.Range("A295").Value = 'برای جمعیتی معادل ۱۰۰۰۰ نفر، نسبت جنسیتی بیش از ۱۰۰۰ به چه معناست؟'
$ws.Rows.Item(295).AutoFit()
$ws.Range("A296").Value = 'هنگام انجام هیسترکتومی از طریق لیگامان پهن، کدام ساختار در معرض آسیب قرار می‌گیرد؟'
$ws.Rows.Item(296).AutoFit()
$ws.Range("A297").Value = 'ویژگی مشخصه تومور لوب فرونتال چیست؟'
$ws.Rows.Item(297).AutoFit()
$ws.Range("A298").Value = 'اکسیژن 100٪ در کدام مورد مؤثر نیست؟ مارس 2005'
$ws.Rows.Item(298).AutoFit()
$ws.Range("A299").Value = 'اسفنکتر تحتانی مری توسط کدام یک از موارد زیر شل می‌شود؟  
الف) الکل  
ب) کافئین / چای  
ج) دیازپام  
د) آنتی‌اسید'
$ws.Rows.Item(299).AutoFit()
$ws.Range("A300").Value = 'در کدام یک از موارد زیر، مقاومت دارویی به دلیل غیرفعال شدن توسط آنزیم‌ها است؟'
$ws.Rows.Item(300).AutoFit()
$ws.Range("A301").Value = 'کدام یک از گزینه‌های زیر به عنوان تجاوز جنسی در نظر گرفته می‌شود؟'
$ws.Rows.Item(301).AutoFit()
$ws.Range("A302").Value = 'اسیدهای چرب توسط تمام بافت های زیر اکسید می شوند، به جز'
$ws.Rows.Item(302).AutoFit()
$ws.Range("A303").Value = 'مواد قابض موادی هستند که:'
$ws.Rows.Item(303).AutoFit()
$ws.Range("A304").Value = 'کدام یک از موارد زیر از طریق gaba عمل نمی‌کند؟'
$ws.Rows.Item(304).AutoFit()
$ws.Range("A305").Value = 'کدام یک از موارد زیر از اثرات تشکیل گاز در فرآیند تجزیه جسد نیست؟'
$ws.Rows.Item(305).AutoFit()
$ws.Range("A306").Value = 'متغیری که می‌توان آن را کنترل یا دستکاری کرد:'
$ws.Rows.Item(306).AutoFit()
$ws.Range("A307").Value = 'در احیای نوزادان، نسبت فشرده‌سازی قفسه سینه به تهویه چیست؟'
$ws.Rows.Item(307).AutoFit()
$ws.Range("A308").Value = 'همه موارد زیر با سندرم تخمدان پلی کیستیک مرتبط هستند، به جز:'
$ws.Rows.Item(308).AutoFit()
$ws.Range("A309").Value = 'یک مرد جوان در طول چند ساعت گذشته درد قفسه سینه رترواسترنال را گزارش می‌دهد. درد به . مرتبط نیست و با نشستن و خم شدن به جلو به طور قابل توجهی تسکین می‌یابد. احتمالاً این حالت ناشی از چیست؟'
$ws.Rows.Item(309).AutoFit()
$ws.Range("A310").Value = 'کندروماتوز سینوویال معمولاً کدام ناحیه را تحت تأثیر قرار می‌دهد؟'
$ws.Rows.Item(310).AutoFit()
$ws.Range("A311").Value = 'دفورمیتی "مداد در فنجان" در کدام بیماری دیده می‌شود؟'
$ws.Rows.Item(311).AutoFit()
$ws.Range("A312").Value = 'تکنیک نیروی متعادل برای آماده‌سازی دندان استفاده می‌شود:'
$ws.Rows.Item(312).AutoFit()
$ws.Range("A313").Value = 'کدام یک از داروهای ضد گلوکوم زیر در نوزادان ناامن است؟'
$ws.Rows.Item(313).AutoFit()
$ws.Range("A314").Value = 'شایع‌ترین تومور خوش‌خیم مری-'
$ws.Rows.Item(314).AutoFit()
$ws.Range("A315").Value = 'این نوار قلب متعلق به یک دانشجوی پزشکی 23 ساله است که از ضربان نامنظم قلب شکایت داشت. به جز نبض نامنظم، قلب او از نظر بالینی طبیعی بود. نوار قلب چه چیزی را نشان می‌دهد؟'
$ws.Rows.Item(315).AutoFit()
$ws.Range("A316").Value = 'فیستول مقعدی بالا یا پایین بر اساس کدام یک از موارد زیر تقسیم‌بندی می‌شود؟'
$ws.Rows.Item(316).AutoFit()
$ws.Range("A317").Value = 'اپیتلیوم انتقالی در همه موارد زیر دیده می‌شود به جز -'
$ws.Rows.Item(317).AutoFit()
$ws.Range("A318").Value = 'چهار استخوان کارپال در چه سنی وجود دارند؟'
$ws.Rows.Item(318).AutoFit()
$ws.Range("A319").Value = 'بیمار الکلی احساس میکند حشرات روی بدنش میخزند، این حالت به عنوان چه شناخته میشود؟'
$ws.Rows.Item(319).AutoFit()
$ws.Range("A320").Value = 'پیش‌آگهی در تومور بافت نرم توسط چه عاملی تعیین می‌شود؟'
$ws.Rows.Item(320).AutoFit()
$ws.Range("A321").Value = 'مقاومت به متوترکسات به چه علتی رخ می‌دهد؟'
$ws.Rows.Item(321).AutoFit()
$ws.Range("A322").Value = 'سندرم شبه آپاندیسیت ناشی از کدام عامل است؟'
$ws.Rows.Item(322).AutoFit()
$ws.Range("A323").Value = 'عبارت صحیح درباره روماتیسمی شدن دریچه میترال -'
$ws.Rows.Item(323).AutoFit()
$ws.Range("A324").Value = 'زنجیره سرد معکوس برای کدام مورد استفاده می‌شود؟'
$ws.Rows.Item(324).AutoFit()
$ws.Range("A325").Value = 'anti dnaase-b برای تشخیص کدام مورد استفاده می‌شود؟'
$ws.Rows.Item(325).AutoFit()
$ws.Range("A326").Value = 'کدام یک از موارد زیر یک انتقال دهنده عصبی مهاری است که به طور انحصاری در سیستم عصبی مرکزی یافت می‌شود؟'
$ws.Rows.Item(326).AutoFit()
$ws.Range("A327").Value = 'بیمار مبتلا به زردی انسدادی ناشی از سرطان پانکراس ممکن است تمام یافته‌های بالینی زیر را داشته باشد به جز'
$ws.Rows.Item(327).AutoFit()
$ws.Range("A328").Value = 'ویژگی شکستگی حفره جمجمه قدامی: wb 07'
$ws.Rows.Item(328).AutoFit()
$ws.Range("A329").Value = 'علامت نیکولادونی برانهام در کدام مورد مشاهده می‌شود؟ مارس 2005'
$ws.Rows.Item(329).AutoFit()
$ws.Range("A330").Value = 'مردی با تومور پانکوست دچار سندرم هورنر شده است. همه موارد زیر از ویژگی‌های سندرم هورنر هستند، به جز:'
$ws.Rows.Item(330).AutoFit()
$ws.Range("A331").Value = 'همه موارد زیر از ویژگی‌های بیماری گودپاسچر هستند به جز-'
$ws.Rows.Item(331).AutoFit()
$ws.Range("A332").Value = 'متیل دوپا در درمان کدام مورد مفید است؟'
$ws.Rows.Item(332).AutoFit()
$ws.Range("A333").Value = 'یک پسر ۹ ساله به مدت ۲ روز دچار کهیرهای خارشدار و قرمز در سراسر بدن شده است. هیچ مشکل تنفسی وجود ندارد. بهترین درمان کدام است؟'
$ws.Rows.Item(333).AutoFit()
$ws.Range("A334").Value = 'ایکتیوزیس ناشی از کدام مورد است؟'
$ws.Rows.Item(334).AutoFit()
$ws.Range("A335").Value = 'مقطع عرضی فایل k-flex چیست؟ (یا) فایل k-flex به طور سنتی از چه جنسی ساخته می‌شود؟'
$ws.Rows.Item(335).AutoFit()
$ws.Range("A336").Value = 'یک کروموزوم حلقوی شکل خاصی از کدام مورد است؟'
$ws.Rows.Item(336).AutoFit()
$ws.Range("A337").Value = 'عرض باکولینگوال دندان های خلفی باید چگونه باشد؟'
$ws.Rows.Item(337).AutoFit()
$ws.Range("A338").Value = 'طبق قانون سلامت روان 2017، حداکثر تعداد روزهایی که یک بیمار می‌تواند به صورت داوطلبانه در بیمارستان روانی بستری شود، چند روز است؟'
$ws.Rows.Item(338).AutoFit()
$ws.Range("A339").Value = 'نادرست در مورد زخم بدخیم معده کدام است؟'
$ws.Rows.Item(339).AutoFit()
$ws.Range("A340").Value = 'کدام دارو از تبدیل t4 به t3 جلوگیری می‌کند؟'
$ws.Rows.Item(340).AutoFit()
$ws.Range("A341").Value = 'مزیت استریلیزاسیون لاپاراسکوپی نسبت به استریلیزاسیون مینی لاپاراتومی چیست؟'
$ws.Rows.Item(341).AutoFit()
$ws.Range("A342").Value = 'وظیفه شغلی دستیاران سلامت مرد چیست؟'
$ws.Rows.Item(342).AutoFit()
$ws.Range("A343").Value = 'کدام یک از موارد زیر باعث ایجاد شبه‌تومور مغزی می‌شود؟'
$ws.Rows.Item(343).AutoFit()
$ws.Range("A344").Value = 'ضایعه کاویته در ریه در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Rows.Item(344).AutoFit()
$ws.Range("A345").Value = 'درست در مورد mhc کلاس ii'
$ws.Rows.Item(345).AutoFit()
$ws.Range("A346").Value = 'همه موارد زیر از عوامل خطر زایمان زودرس هستند به جز:'
$ws.Rows.Item(346).AutoFit()
$ws.Range("A347").Value = 'یک پسر 27 ساله به مدت 3 سال دچار پاپول‌های خارش‌دار و خراشیده روی پیشانی و قسمت‌های در معرض بازوها و پاها شده است. بیماری در فصل بارانی شدیدتر بوده و در زمستان کاملاً بهبود می‌یابد. محتمل‌ترین تشخیص چیست؟'
$ws.Rows.Item(347).AutoFit()
$ws.Range("A348").Value = 'کوتاهی متاکارپال چهارم از ویژگی‌های کدام یک از موارد زیر است؟'
$ws.Rows.Item(348).AutoFit()
$ws.Range("A349").Value = 'مدت زمان بی‌حسی نخاعی به همه موارد زیر بستگی دارد به جز:'
$ws.Rows.Item(349).AutoFit()
$ws.Range("A350").Value = 'در داده‌های جمع‌آوری شده از روش نیتروژن شست‌وشو در یک بیمار، حجم نیتروژن در حجم باقی‌مانده (rv) او 800 میلی‌لیتر بود. حجم باقی‌مانده (rv) در این بیمار چقدر است؟'
$ws.Rows.Item(350).AutoFit()
$ws.Range("A351").Value = 'کدام یک از موارد زیر را می‌توان به عنوان انواع مختلف مرگ در نظر گرفت؟'
$ws.Rows.Item(351).AutoFit()
$ws.Range("A352").Value = 'کدام یک از موارد زیر نشانه بستری شدن در پنومونی نیست؟'
$ws.Rows.Item(352).AutoFit()
$ws.Range("A353").Value = 'درمان دارویی انتخابی برای sspe کدام است؟'
$ws.Rows.Item(353).AutoFit()
$ws.Range("A354").Value = 'کدام تاندون برای ترمیم رباط صلیبی قدامی استفاده می‌شود؟'
$ws.Rows.Item(354).AutoFit()
$ws.Range("A355").Value = 'تغییر رنگ زرد پوست و مخاط در کدام مسمومیت دیده می‌شود؟'
$ws.Rows.Item(355).AutoFit()
$ws.Range("A356").Value = 'چه چیزی از دندان نهفته خارج می‌شود؟'
$ws.Rows.Item(356).AutoFit()
$ws.Range("A357").Value = 'عبارت صحیح درباره تیتر آنتی استرپتولیزین o کدام است؟'
$ws.Rows.Item(357).AutoFit()
$ws.Range("A358").Value = 'کدام آزمایش برای تشخیص ایسکمی قابل برگشت میوکارد انجام می‌شود؟'
$ws.Rows.Item(358).AutoFit()
$ws.Range("A359").Value = 'یک زن ۶۲ ساله تحت عمل پانکراتیکودودنکتومی برای سرطان سر پانکراس قرار می‌گیرد. یک ژژنوستومی برای تسهیل تغذیه جایگذاری می‌شود زیرا انتظار می‌رود بهبودی او طولانی باشد. بهترین روش برای ارائه تغذیه پس از عمل چیست؟'
$ws.Rows.Item(359).AutoFit()
$ws.Range("A360").Value = 'دوره نهفتگی متوسط برای شانکر سیفلیس اولیه چند است؟'
$ws.Rows.Item(360).AutoFit()
$ws.Range("A361").Value = 'همه موارد زیر در مورد مری بارت صحیح هستند به جز'
$ws.Rows.Item(361).AutoFit()
$ws.Range("A362").Value = 'کارسینومای سلول کلیوی دوطرفه در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Rows.Item(362).AutoFit()
$ws.Range("A363").Value = 'همه موارد زیر در مورد آمیلوئیدوز صحیح هستند، به جز:'
$ws.Rows.Item(363).AutoFit()
$ws.Range("A364").Value = 'شکستگی بنت، دررفتگی شکستگی قاعده کدام متاکارپال است:'
$ws.Rows.Item(364).AutoFit()
$ws.Range("A365").Value = 'یک زن جوان با سابقه تنگی نفس هنگام فعالیت مراجعه می‌کند. در معاینه، شکاف ثابت و گسترده s2 با سوفل سیستولیک اِجکشن (iii/vi) در فضای بین دنده‌ای دوم چپ مشاهده می‌شود. نوار قلب او انحراف محور به چپ را نشان می‌دهد. محتمل‌ترین تشخیص چیست؟'
$ws.Rows.Item(365).AutoFit()
$ws.Range("A366").Value = 'سه‌گانه سامتر به ارتباط بین آسم حساس به آسپرین و کدام یک از موارد زیر اشاره دارد؟'
$ws.Rows.Item(366).AutoFit()
$ws.Range("A367").Value = 'در معاینه دهان یک مرد 57 ساله، یک ناحیه بزرگ کراتوتیک که تمام کام را پوشانده است مشاهده می‌شود. همچنین برخی نقاط قرمز در این ناحیه دیده می‌شود. به احتمال زیاد این بیمار:'
$ws.Rows.Item(367).AutoFit()
$ws.Range("A368").Value = 'گنوکوک‌ها ابتدا کدام ناحیه را آلوده می‌کنند؟'
$ws.Rows.Item(368).AutoFit()
$ws.Range("A369").Value = 'آدنوکارسینوم سینوس اتموئید با چه گروهی از کارگران مرتبط است؟'
$ws.Rows.Item(369).AutoFit()
$ws.Range("A370").Value = 'مدروکسی پروژسترون استات دپو چیست؟'
$ws.Rows.Item(370).AutoFit()
$ws.Range("A371").Value = 'مننژیت مننگوکوکی در کمبود کدام یک از اجزای کمپلمان مشاهده می‌شود؟'
$ws.Rows.Item(371).AutoFit()
$ws.Range("A372").Value = 'کدام یک از دیورتیک‌های زیر مهارکننده سمپورتر na+.k+.2cl- نیست؟'
$ws.Rows.Item(372).AutoFit()
$ws.Range("A373").Value = 'افزایش شکنندگی اسمزی در کدام یک از موارد زیر مشاهده می‌شود؟'
$ws.Rows.Item(373).AutoFit()
$ws.Range("A374").Value = 'هیستامین همه موارد زیر را باعث می‌شود به جز -'
$ws.Rows.Item(374).AutoFit()
$ws.Range("A375").Value = 'مهم‌ترین منبع خونرسانی به معده کدام است؟'
$ws.Rows.Item(375).AutoFit()
$ws.Range("A376").Value = 'درمان انتخابی چیست؟'
$ws.Rows.Item(376).AutoFit()
$ws.Range("A377").Value = 'پدر نظارت اپیدمیولوژیک مدرن کیست؟'
$ws.Rows.Item(377).AutoFit()
$ws.Range("A378").Value = 'مری در کدام سطح شروع می‌شود؟'
$ws.Rows.Item(378).AutoFit()
$ws.Range("A379").Value = 'مارپیچ‌های کورشمن در خلط در کدام مورد دیده می‌شوند؟'
$ws.Rows.Item(379).AutoFit()
$ws.Range("A380").Value = 'کدام یک از موارد زیر دیورتیکول کاذب هستند؟'
$ws.Rows.Item(380).AutoFit()
$ws.Range("A381").Value = 'سندرم دامپینگ به دلیل کدام مورد رخ می‌دهد؟'
$ws.Rows.Item(381).AutoFit()
$ws.Range("A382").Value = 'تعریف کودک autistic (اوتیستیک) چیست؟'
$ws.Rows.Item(382).AutoFit()
$ws.Range("A383").Value = 'کدام یک از داروهای زیر منجر به عارضه ذکر شده می‌شود؟'
$ws.Rows.Item(383).AutoFit()
$ws.Range("A384").Value = 'کدام یک از شریان های زیر مشتق از قوس آئورتی دوم است؟'
$ws.Rows.Item(384).AutoFit()
$ws.Range("A385").Value = 'شایع‌ترین محل لنفوم در دستگاه گوارش کدام است؟'
$ws.Rows.Item(385).AutoFit()
$ws.Range("A386").Value = 'عصب‌دهی حسی حنجره، زیر تارهای صوتی -'
$ws.Rows.Item(386).AutoFit()
$ws.Range("A387").Value = 'میگزومای قلبی بیشتر در کدام قسمت دیده می‌شود؟'
$ws.Rows.Item(387).AutoFit()
$ws.Range("A388").Value = 'یک خانم با توده شکمی مورد بررسی قرار گرفت. توده های تخمدانی دو طرفه با سطح صاف در او یافت شد. در بررسی میکروسکوپی، سلول های با شکل حلقه نگین دار مشاهده شد. تشخیص چیست؟'
$ws.Rows.Item(388).AutoFit()
$ws.Range("A389").Value = 'رشد مو در تمامی این نواحی با سطح آندروژن سرم مرتبط است به جز...'
$ws.Rows.Item(389).AutoFit()
$ws.Range("A390").Value = 'تمامی عبارات زیر در مورد سدیم فلوراید در اتواسکلروز صحیح هستند، به جز'
$ws.Rows.Item(390).AutoFit()
$ws.Range("A391").Value = 'کدام یک از عوامل بیهوشی زیر می‌تواند باعث تاکیپنه شود؟'
$ws.Rows.Item(391).AutoFit()
$ws.Range("A392").Value = 'شایع‌ترین محل کارسینومای زبان کدام است؟'
$ws.Rows.Item(392).AutoFit()
$ws.Range("A393").Value = 'در یک زن با ایمپلنت زیرجلدی پروژسترون، اختلال قاعدگی مشاهده شده چیست؟'
$ws.Rows.Item(393).AutoFit()
$ws.Range("A394").Value = 'نشانگر واسکولیت کلیوی در کودکان کدام است؟'
$ws.Rows.Item(394).AutoFit()
$ws.Range("A395").Value = 'کدام یک از موارد زیر توسط گنوکوک تحت تأثیر قرار نمی‌گیرد؟'
$ws.Rows.Item(395).AutoFit()
$ws.Range("A396").Value = 'دفورمیتی گردن قو در کدام یک از موارد زیر مشاهده می‌شود؟'
$ws.Rows.Item(396).AutoFit()
$ws.Range("A397").Value = 'پنوماتوسل‌ها در رادیوگرافی قفسه سینه در یک کودک مبتلا به پنومونی در عفونت با کدام عامل دیده می‌شوند؟'
$ws.Rows.Item(397).AutoFit()
$ws.Range("A398").Value = 'کدام مورد در مجموعه آزمون‌های هالستید-ریتان وجود ندارد؟('
$ws.Rows.Item(398).AutoFit()
$ws.Range("A399").Value = 'لایسین در کدام ماده کمبود دارد؟'
$ws.Rows.Item(399).AutoFit()
$ws.Range("A400").Value = 'شایع ترین غده بزاقی که دچار سنگ می شود:'
$ws.Rows.Item(400).AutoFit()
$ws.Range("A401").Value = 'تست کومبس مستقیم در کم خونی همولیتیک به دلیل کدام مورد مثبت می‌شود؟'
$ws.Rows.Item(401).AutoFit()
$ws.Range("A402").Value = 'یک زن ۷۵ ساله با سابقه فشار خون بالا، دچار خستگی و تنگی نفس هنگام فعالیت شده است. فشار خون او ۱۶۰/۶۰ میلی‌متر جیوه و نبض ۸۰ در دقیقه است. صدای دوم قلب کاهش یافته و یک سوفل دیاستولیک اولیه دارد که از کنار راست جناغ به سمت نوک قلب انتشار می‌یابد. تشخیص بالینی شما نارسایی آئورت است. برای این بیمار، یافته مشخصه نبض شریانی را انتخاب کنید.'
$ws.Rows.Item(402).AutoFit()
$ws.Range("A403").Value = 'درمان انتخابی بیماری کاوازاکی کدام است؟'
$ws.Rows.Item(403).AutoFit()
$ws.Range("A404").Value = 'فتق غیر مرتبط با دیواره شکم:-'
$ws.Rows.Item(404).AutoFit()
$ws.Range("A405").Value = 'مسمومیت با سولفات منیزیم شامل همه موارد زیر است به جز'
$ws.Rows.Item(405).AutoFit()
$ws.Range("A406").Value = 'همه موارد زیر جزو نوروتیک‌های مغزی هستند به جز'
$ws.Rows.Item(406).AutoFit()
$ws.Range("A407").Value = 'هایپرآلدوسترونیسم با موارد زیر مشخص می‌شود به جز -'
$ws.Rows.Item(407).AutoFit()
$ws.Range("A408").Value = 'شایع‌ترین محل گاسترینوم کجاست؟'
$ws.Rows.Item(408).AutoFit()
$ws.Range("A409").Value = 'نوع نفروپاتی در دیابت ملیتوس-'
$ws.Rows.Item(409).AutoFit()
$ws.Range("A410").Value = 'نادرست در مورد همانژیوم کاورنوتی کدام است؟'
$ws.Rows.Item(410).AutoFit()
$ws.Range("A411").Value = 'یک دانشجوی سال اول کالج 18 ساله اتاق خوابگاه خود را با سه هم اتاقی تقسیم می‌کند. یکی از هم‌اتاقی‌های او به مننژیت مننگوکوکی، ناشی از باکتری نایسریا مننژیتیدیس، مبتلا شده است. سه هم‌اتاقی دیگر قرنطینه شده و روزانه دو بار با یک آنتی‌بیوتیک به عنوان پیشگیری در برابر این باکتری تحت درمان قرار می‌گیرند، زیرا هیچ‌یک از آنها پیش از ثبت‌نام واکسن مننگوکوک را دریافت نکرده بودند. به آنها گفته می‌شود که این آنتی‌بیوتیک می‌تواند باعث تغییر رنگ قرمز در ادرار یا اشک آنها شود. دلیل مؤثر بودن این دارو در کشتن باکتری کدام یک از موارد زیر است؟'
$ws.Rows.Item(411).AutoFit()
$ws.Range("A412").Value = 'همه موارد زیر از علائم قرنیه در تراخم هستند به جز-'
$ws.Rows.Item(412).AutoFit()
$ws.Range("A413").Value = 'فلج دوطرفه عصب حنجره‌ای بازگشتی در کدام مورد دیده می‌شود؟'
$ws.Rows.Item(413).AutoFit()
$ws.Range("A414").Value = 'بینایی تونلی ناشی از کدام مورد است؟'
$ws.Rows.Item(414).AutoFit()
$ws.Range("A415").Value = 'کدام یک از داروهای زیر در درمان یووئیت قدامی حاد استفاده نمی‌شود؟'
$ws.Rows.Item(415).AutoFit()
$ws.Range("A416").Value = 'کنژنکتیویت لیگنئوس ناشی از چیست؟'
$ws.Rows.Item(416).AutoFit()
$ws.Range("A417").Value = 'پنوموسیستیس کارینی یک قارچ است زیرا؟'
$ws.Rows.Item(417).AutoFit()
$ws.Range("A418").Value = 'چرا واکسن های زیر رنگی هستند؟'
$ws.Rows.Item(418).AutoFit()
$ws.Range("A419").Value = 'علت سقط جنین مکرر در سه‌ماهه دوم بارداری:'
$ws.Rows.Item(419).AutoFit()
$ws.Range("A420").Value = 'نقش ماری‌جوانا در کاشکسی مرتبط با ایدز چیست؟'
$ws.Rows.Item(420).AutoFit()
$ws.Range("A421").Value = 'ph واژن در چه زمانی در پایین‌ترین حد خود است؟'
$ws.Rows.Item(421).AutoFit()
$ws.Range("A422").Value = 'سوفل پیوسته در تمام موارد زیر یافت می‌شود، به جز -'
$ws.Rows.Item(422).AutoFit()
$ws.Range("A423").Value = 'معیارهای غیرقابل برداشت بودن در سرطان مجرای صفراوی هیلار همه موارد زیر هستند به جز'
$ws.Rows.Item(423).AutoFit()
$ws.Range("A424").Value = 'نفروپاتی ریفلاکس مزمن باعث کدام یک از موارد زیر می‌شود؟'
$ws.Rows.Item(424).AutoFit()
$ws.Range("A425").Value = 'ماده حاجب انتخابی برای میلوگرام کدام است؟'
$ws.Rows.Item(425).AutoFit()
$ws.Range("A426").Value = '''پارگی دوازدهه'' چیست؟'
$ws.Rows.Item(426).AutoFit()
$ws.Range("A427").Value = 'غلاف میلین در سیستم عصبی مرکزی توسط چه چیزی تولید می‌شود؟'
$ws.Rows.Item(427).AutoFit()
$ws.Range("A428").Value = 'افزایش کدام یک از پارامترهای زیر منحنی تفکیک اکسیژن را به چپ منتقل می‌کند؟'
$ws.Rows.Item(428).AutoFit()
$ws.Range("A429").Value = 'ویژگی‌های سندرم پوتز-جگرز شامل همه موارد زیر است به جز'
$ws.Rows.Item(429).AutoFit()
$ws.Range("A430").Value = 'مرواریدهای عنبیه در کدام بیماری دیده می‌شوند؟'
$ws.Rows.Item(430).AutoFit()
$ws.Range("A431").Value = 'تنه تیروسرویکال شاخه ای از کدام قسمت است؟'
$ws.Rows.Item(431).AutoFit()
$ws.Range("A432").Value = 'کدام یک از عوامل فیبرینولیتیک زیر آنتی ژنیک است؟'
$ws.Rows.Item(432).AutoFit()
$ws.Range("A433").Value = 'یک زن باردار در هفته ۳۷ بارداری به یک مرکز بهداشت اولیه برای معاینه قبل از زایمان مراجعه می‌کند. او تاکنون هیچ معاینه قبل از زایمانی نداشته است. بهترین اقدام در مورد ایمن‌سازی کزاز در این مورد چیست؟'
$ws.Rows.Item(433).AutoFit()
$ws.Range("A434").Value = 'کدام واکسن به حرارت حساس‌تر است؟'
$ws.Rows.Item(434).AutoFit()
$ws.Range("A435").Value = 'روش تشخیصی انتخابی برای پاتولوژی پاراتیروئید کدام است؟'
$ws.Rows.Item(435).AutoFit()
$ws.Range("A436").Value = 'داروی انتخابی برای اکتینومایکوز قفسه سینه کدام است؟'
$ws.Rows.Item(436).AutoFit()
$ws.Range("A437").Value = 'در مورد دراکونکولیازیس صحیح نیست کدام گزینه؟'
$ws.Rows.Item(437).AutoFit()
$ws.Range("A438").Value = 'غنی‌ترین منبع اسیدهای چرب ضروری'
$ws.Rows.Item(438).AutoFit()
$ws.Range("A439").Value = 'یک پسر ۱۶ ساله توسط مربی مدرسه برای معاینه فیزیکی قبل از پیوستن به تیم فوتبال ارجاع داده شده است. برادر بزرگتر او به طور ناگهانی در حین تمرین فوتبال فوت کرده است، کالبد شکافی انجام نشده است. بیمار در معاینه قفسه سینه یک سوفل سیستولیک بلند دارد. همه موارد زیر با کاردیومیوپاتی هیپرتروفیک سازگار هستند، به جز'
$ws.Rows.Item(439).AutoFit()
$ws.Range("A440").Value = 'منشأ مجرای تیروگلوسال کدام است؟'
$ws.Rows.Item(440).AutoFit()
$ws.Range("A441").Value = 'یک شریانچه با لایه سلول‌های اندوتلیال آسیب‌دیده قادر به انجام کدام مورد نخواهد بود؟'
$ws.Rows.Item(441).AutoFit()
$ws.Range("A442").Value = 'امپیما اصطلاحی است که برای چه چیزی به کار می‌رود؟'
$ws.Rows.Item(442).AutoFit()
$ws.Range("A443").Value = 'شایع ترین ویروس مولوسکوم:'
$ws.Rows.Item(443).AutoFit()
$ws.Range("A444").Value = 'پراکسیداسیون لیپیدهای چند غیراشباع غشاهای زیرسلولی چه چیزی تولید می‌کند؟'
$ws.Rows.Item(444).AutoFit()
$ws.Range("A445").Value = 'کدام یک از موارد زیر از موارد منع زایمان طبیعی پس از سزارین قبلی نیست؟'
$ws.Rows.Item(445).AutoFit()
$ws.Range("A446").Value = 'تعداد روزهای مورد نیاز برای تشخیص هیپومانیا عبارت است از:'
$ws.Rows.Item(446).AutoFit()
$ws.Range("A447").Value = 'نمای‌های استاندارد برای ماموگرافی کدام‌اند؟'
$ws.Rows.Item(447).AutoFit()
$ws.Range("A448").Value = 'یک بیمار 2 ماه پیش عمل شده و در آن زمان از برش میدلاین استفاده شده است. اکنون به عمل دوم نیاز دارد و بهترین برش برای استفاده در این زمان کدام است؟'
$ws.Rows.Item(448).AutoFit()
$ws.Range("A449").Value = 'یک کودک با نارسایی گردش خون محیطی مراجعه می‌کند. ph شریانی 7.0، pco2 برابر 15 mmhg و po2 برابر 76 mmhg است. کدام یک از موارد زیر درمان فوری خواهد بود؟'
$ws.Rows.Item(449).AutoFit()
$ws.Range("A450").Value = 'پیشگیری اولیه در سکته قلبی شامل همه موارد زیر به جز کدام است؟'
$ws.Rows.Item(450).AutoFit()
$ws.Range("A451").Value = 'درصد مواد معدنی در مینای کاملاً رشد یافته حدوداً چقدر است؟'
$ws.Rows.Item(451).AutoFit()
$ws.Range("A452").Value = 'بیمار با سردرد و گرگرفتگی مراجعه کرده است. وی سابقه خانوادگی فوت یکی از بستگانش به دلیل تومور تیروئید را دارد. بررسی مورد نیاز برای این بیمار کدام است؟'
$ws.Rows.Item(452).AutoFit()
$ws.Range("A453").Value = 'خطر کارسینوم در کولیت اولسراتیو:'
$ws.Rows.Item(453).AutoFit()
$ws.Range("A454").Value = 'فشار خون بدخیم باعث کدام یک از موارد زیر می‌شود؟'
$ws.Rows.Item(454).AutoFit()
$ws.Range("A455").Value = 'ویتامین k چه اسید آمینه ای را کربوکسیله می کند؟'
$ws.Rows.Item(455).AutoFit()
$ws.Range("A456").Value = 'afp در کدام یک از موارد زیر افزایش می‌یابد؟  
الف) تومور کیسه زرده  
ب) سمینوما  
ج) تراتوما  
د) تومور سینوس اندودرمال  
ه) سیروز'
$ws.Rows.Item(456).AutoFit()
$ws.Range("A457").Value = 'مقدار قطعی طول دهانه رحم برای پیش‌بینی زایمان زودرس چیست؟'
$ws.Rows.Item(457).AutoFit()
$ws.Range("A458").Value = 'پاسخ مردمک همیانوپیک ورنیکه در ضایعات کدام ناحیه دیده می‌شود؟'
$ws.Rows.Item(458).AutoFit()
$ws.Range("A459").Value = 'فرآیند بیماری که به بهترین شکل این مشکل را توضیح می‌دهد:'
$ws.Rows.Item(459).AutoFit()
$ws.Range("A460").Value = 'یک زن 40 ساله با cin 2، مرحله بعدی مدیریت چیست؟'
$ws.Rows.Item(460).AutoFit()
$ws.Range("A461").Value = 'سلول‌های آماکرین چه ماده‌ای ترشح می‌کنند؟'
$ws.Rows.Item(461).AutoFit()
$ws.Range("A462").Value = 'اصطلاح ترجمه نیک (nick translation) به چه چیزی اشاره دارد؟'
$ws.Rows.Item(462).AutoFit()
$ws.Range("A463").Value = 'قوی ترین پیوند از میان گزینه های زیر کدام است؟'
$ws.Rows.Item(463).AutoFit()
$ws.Range("A464").Value = 'استخوانی شدن فونتانل قدامی در چه زمانی اتفاق می‌افتد؟ مارس 2007'
$ws.Rows.Item(464).AutoFit()
$ws.Range("A465").Value = 'یک پسر 15 ساله با سابقه یک روزه خونریزی لثه، خونریزی زیر ملتحمه و راش پورپوریک مراجعه کرده است. نتایج بررسی‌ها به شرح زیر است: هموگلوبین 6.4 گرم در دسی‌لیتر، گلبول‌های سفید 26500 در میلی‌متر مکعب، پلاکت 35000 در میلی‌متر مکعب، زمان پروترومبین 20 ثانیه با کنترل 13 ثانیه، زمان ترومبوپلاستین جزئی 50 ثانیه و فیبرینوژن 10 میلی‌گرم در دسی‌لیتر. اسمیر محیطی نشان‌دهنده لوسمی میلوبلاستیک حاد بود. کدام یک از گزینه‌های زیر محتمل‌ترین تشخیص است؟'
$ws.Rows.Item(465).AutoFit()
$ws.Range("A466").Value = 'کدام یک از سرطان‌های زیر به طور شایع به کبد متاستاز نمی‌دهد؟'
$ws.Rows.Item(466).AutoFit()
$ws.Range("A467").Value = 'رنگ سیلندر نیتروس اکساید چیست؟'
$ws.Rows.Item(467).AutoFit()
